$p = $ppt.ActivePresentation

# Locate the shape/paragraph that contains the literal text "(lattice version)".
# (Known to be "CustomShape 14" on slide 1, but search defensively instead of
# hard-coding indices in case shape ordering ever differs.)
$targetShape = $null
$targetSlide = $null
$targetParaIndex = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $paraCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $paraCount; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    if ($para.Text -like "*(lattice version)*") {
                        $targetShape = $shape
                        $targetSlide = $slide
                        $targetParaIndex = $pi
                    }
                }
            }
        }
    }
}

$tf = $targetShape.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs($targetParaIndex, 1)

# Split the run right after the 'l' in "(lattice version)" by inserting and
# then removing a one-character marker -- this leaves a genuine run boundary
# (mirroring what PowerPoint itself does when a run is edited mid-word)
# without disturbing the run formatting (font, size, bold, color, ...).
$afterL = $para.Characters(2, 1)
[void]$afterL.InsertAfter("|")
$marker = $para.Characters(3, 1)
$marker.Text = ""

# Drop the surrounding parentheses so the text reads "lattice version".
$lastChar = $para.Characters(17, 1)
$lastChar.Text = ""
$firstChar = $para.Characters(1, 1)
$firstChar.Text = ""
